# DataSource - Emision Motor NuevasCoberturas.xlsx
# "se modificad data para regresion en preprod R33"
#
# Update the regression-test data row on the active sheet (Sheet1):
#  - Ambiente/URL columns switch from the "preproducciongestion" host
#    to the "i-preproducciongestion" (pre-prod) host
#  - NroCuenta test value changes
#  - FechaInicio test value changes
#  - Selection moves from E2 to B3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("E2").Value = 5069929970
$ws.Range("K2").Value = "17/06/2021"

$ws.Range("B3").Select()
